# Generate Report for Handback
# The "fe5756af-a198-49bf-80ce-74ebceea6cc5.md" file has now been handed
# back (in sync with en-US) for both the zh-cn and de-de locales, so the
# status + handback timestamps on the report need to be refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the fe5756af file; both locale status columns
# move from "Ready for handoff" to "Handed back: in sync with en-US".
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn detail sheet: row 3 is the fe5756af file. Status + handback
# datetime are refreshed.
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-18 02:57:19"

# de-de detail sheet: row 3 is the fe5756af file. Status + handback
# datetime are refreshed.
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-18 02:57:31"
